# Generate Report for Handoff
# Regenerates the handoff report: the four "Ready for handoff" files move
# from "low" to "ht" priority and get a fresh handoff timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the four files that were
# just handed off again.
$overview.Range("G4:G7").Value = "2016-08-15 08:47:31"

# zh-cn sheet: Priority + Latest Handoff Datetime refreshed for rows 4-7.
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-08-15 08:47:26"
$zhcn.Range("H5").Value = "2016-08-15 08:47:26"
$zhcn.Range("H6").Value = "2016-08-15 08:47:26"
$zhcn.Range("H7").Value = "2016-08-15 08:47:26"

# de-de sheet: Priority refreshed for rows 4-7. The "Latest Handoff Datetime"
# text for these rows shares its underlying value with the Overview sheet's
# "Latest HO Xliff Generate Date" column, so it moves to the same new
# timestamp (08:47:31) rather than the zh-cn sheet's 08:47:26.
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

$dede.Range("H4").Value = "2016-08-15 08:47:31"
$dede.Range("H5").Value = "2016-08-15 08:47:31"
$dede.Range("H6").Value = "2016-08-15 08:47:31"
$dede.Range("H7").Value = "2016-08-15 08:47:31"
